# "update front and back" - challenges.xlsx
#
# The source lake list was rekeyed: the old lowercase "geojson_key" text
# values (genfersee / bielersee / murtensee / neuenburgersee) are replaced
# by the same human-readable names used in column A, the per-row meter
# totals were refreshed, and the Murtensee row was dropped entirely
# (Neuenburgersee moves up into row 4). Column C was also widened and the
# B2 total gets a smaller/darker wrapped font, and the sheet was set up
# for portrait-A4 printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Genfersee ---
$ws.Range("A2").Value = "Genferseeumrundung"
$ws.Range("B2").Value = 190980
$ws.Range("C2").Value = "Genferseeumrundung"

# --- Row 3: Bielersee ---
$ws.Range("A3").Value = "Bielerseeumrundung"
$ws.Range("B3").Value = 516500
$ws.Range("C3").Value = "Bielerseeumrundung"

# --- Row 4: Neuenburgersee (Murtensee row below it gets dropped) ---
$ws.Range("A4").Value = "Neuenburgerseeumrundung"
$ws.Range("B4").Value = 983300
$ws.Range("C4").Value = "Neuenburgerseeumrundung"

# --- Drop the old Murtensee row (was row 5) entirely ---
$ws.Range("A5:C5").EntireRow.Delete()

# --- Column C is widened to fit the longer geojson_key text ---
$ws.Columns("C").ColumnWidth = 41.619791666666664

# --- B2 (Genfersee total_meter) gets a smaller, dark, wrapped font and loses its border ---
$ws.Range("B2").Font.Size = 9
$ws.Range("B2").Font.Color = 1052688
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").WrapText = $true
$ws.Range("B2").Borders.LineStyle = 0

# --- Page setup: portrait, A4 ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection cosmetically moved ---
[void]$ws.Range("C13").Select()
